# Apply the Jan 31 2024 cryptos-list refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text storage (Quote-Prefix), matching the original
# inlineStr cells and preventing Excel from auto-converting e.g. '99.30' to 99.3.
$q = "'"

$ws.Range('D2').Value = $q + '42.743.49'
$ws.Range('E2').Value = $q + '  -1.04%  '
$ws.Range('D3').Value = $q + '2.309.56'
$ws.Range('E3').Value = $q + '  +0.32%  '
$ws.Range('E4').Value = $q + '  -0.01%  '
$ws.Range('D5').Value = $q + '302.98'
$ws.Range('E5').Value = $q + '  -1.84%  '
$ws.Range('D6').Value = $q + '99.30'
$ws.Range('E6').Value = $q + '  -4.56%  '
$ws.Range('D7').Value = $q + '0.504'
$ws.Range('E7').Value = $q + '  -4.26%  '
$ws.Range('E8').Value = $q + '  +0.03%  '
$ws.Range('D9').Value = $q + '0.502'
$ws.Range('E9').Value = $q + '  -2.92%  '
$ws.Range('D10').Value = $q + '34.68'
$ws.Range('E10').Value = $q + '  -2.83%  '
$ws.Range('D11').Value = $q + '0.0792'
$ws.Range('E11').Value = $q + '  -1.99%  '
$ws.Range('E12').Value = $q + '  +0.59%  '
$ws.Range('D13').Value = $q + '6.70'
$ws.Range('E13').Value = $q + '  -3.27%  '
$ws.Range('D14').Value = $q + '2.669.88'
$ws.Range('E14').Value = $q + '  +0.22%  '
$ws.Range('D15').Value = $q + '15.65'
$ws.Range('D16').Value = $q + '2.304.15'
$ws.Range('E16').Value = $q + '  +0.09%  '
$ws.Range('D17').Value = $q + '0.800'
$ws.Range('E17').Value = $q + '  +0.33%  '
$ws.Range('D18').Value = $q + '42.697.53'
$ws.Range('D19').Value = $q + '0.0₃0905'
$ws.Range('E19').Value = $q + '  -1.51%  '
$ws.Range('E20').Value = $q + '  -3.32%  '
$ws.Range('D21').Value = $q + '6.06'
$ws.Range('E21').Value = $q + '  -1.66%  '
$ws.Range('D22').Value = $q + '67.84'
$ws.Range('E22').Value = $q + '  +0.21%  '
$ws.Range('D23').Value = $q + '234.82'
$ws.Range('E23').Value = $q + '  -1.97%  '
$ws.Range('D24').Value = $q + '1.95'
$ws.Range('E24').Value = $q + '  -2.71%  '
$ws.Range('D25').Value = $q + '2.51'
$ws.Range('E25').Value = $q + '  -3.15%  '
$ws.Range('E26').Value = $q + '  +0.09%  '
$ws.Range('D27').Value = $q + '24.90'
$ws.Range('E27').Value = $q + '  +1.09%  '
$ws.Range('E28').Value = $q + '  -6.74%  '
$ws.Range('D29').Value = $q + '34.47'
$ws.Range('E29').Value = $q + '  -4.17%  '
$ws.Range('D30').Value = $q + '164.23'
$ws.Range('E30').Value = $q + '  +1.52%  '
$ws.Range('D31').Value = $q + '9.12'
$ws.Range('E31').Value = $q + '  -4.34%  '
$ws.Range('D32').Value = $q + '1.00'
$ws.Range('E32').Value = $q + '  +0.02%  '
$ws.Range('D33').Value = $q + '5.00'
$ws.Range('E33').Value = $q + '  -4.29%  '
$ws.Range('E34').Value = $q + '  -4.81%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = $q + '16.81'
$ws.Range('E35').Value = $q + '  -7.60%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = $q + '4.43'
$ws.Range('E36').Value = $q + '  -0.43%  '
$ws.Range('D37').Value = $q + '0.0700'
$ws.Range('E37').Value = $q + '  -4.28%  '
$ws.Range('D38').Value = $q + '2.88'
$ws.Range('E38').Value = $q + '  -3.40%  '
$ws.Range('E39').Value = $q + '  -2.85%  '
$ws.Range('E40').Value = $q + '  -5.20%  '
$ws.Range('D41').Value = $q + '0.110'
$ws.Range('E41').Value = $q + '  -3.64%  '
$ws.Range('D42').Value = $q + '2.50'
$ws.Range('E42').Value = $q + '  -1.20%  '
$ws.Range('D43').Value = $q + '1.967.21'
$ws.Range('E43').Value = $q + '  +0.34%  '
$ws.Range('D44').Value = $q + '0.0279'
$ws.Range('E44').Value = $q + '  -3.06%  '
$ws.Range('D45').Value = $q + '18.52'
$ws.Range('E45').Value = $q + '  -1.07%  '
$ws.Range('D46').Value = $q + '10.19'
$ws.Range('E46').Value = $q + '  +0.49%  '
$ws.Range('D47').Value = $q + '2.87'
$ws.Range('E47').Value = $q + '  -6.07%  '
$ws.Range('D48').Value = $q + '55.21'
$ws.Range('E48').Value = $q + '  -2.77%  '
$ws.Range('D49').Value = $q + '2.535.81'
$ws.Range('E49').Value = $q + '  +0.26%  '
$ws.Range('D50').Value = $q + '2.83'
$ws.Range('E50').Value = $q + '  -3.16%  '
$ws.Range('E51').Value = $q + '  +1.06%  '
